$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 33

$ws.Range("E4").Value = 31
$ws.Range("F4").Value = 14
$ws.Range("H4").Value = 17

$ws.Range("E12").Value = 39
$ws.Range("F12").Value = 15
$ws.Range("H12").Value = 17

$ws.Range("E15").Value = 118

$ws.Range("E16").Value = 339

$ws.Range("E18").Value = 104
